# Applies the "Fix minor formatting issue" edit to the Lending Club data
# story document.
#
# Summary of the change (see commit diff):
#  - Paragraph 1 + 2 (the two "what are some insights.../By now you've
#    asked..." prompt paragraphs) collapse into a single new title
#    paragraph "The Lending Club: Initial Research".
#  - The two now-empty spacer paragraphs are removed.
#  - Every remaining paragraph/run has its direct formatting normalized
#    down to just the Helvetica font (no explicit color/size/shading).
#  - The inline OLE chart object is resized from 427x296pt to
#    454.5x314.5pt (and its (cosmetic) ObjectID bumped).
#  - A couple of runs get split at the same character offsets they were
#    split at upstream (no text changes), and the stray
#    w:lastRenderedPageBreak marker is dropped.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function RunXml([string]$text, [bool]$preserve) {
    $space = ''
    if ($preserve) { $space = ' xml:space="preserve"' }
    return '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/></w:rPr><w:t' + $space + '>' + $text + '</w:t></w:r>'
}

$pPrHelv = '<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/></w:rPr></w:pPr>'

# ---------------------------------------------------------------------
# Work from the bottom of the document upwards so the paragraph indices
# of the not-yet-processed (earlier) paragraphs stay valid while we
# delete / rewrite content further down.
# ---------------------------------------------------------------------

# Paragraph 10: "One possibility is ... " (split around the _GoBack bookmark)
$p10 = $d.Paragraphs.Item(10).Range
$xml10 = '<w:p ' + $wNs + '>' + $pPrHelv +
    (RunXml 'One possibility is that these differences in default rates between states is simply not statistically significant. The other is that the difference is taken into account when assigning ratings to the ' $true) +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    (RunXml 'loans. These possibilities are avenues for further exploration.' $false) +
    '</w:p>'
$p10.InsertXML($xml10)

# Paragraph 9: "Note that the results ... " (split into 4 runs, drop the
# lastRenderedPageBreak marker)
$p9 = $d.Paragraphs.Item(9).Range
$xml9 = '<w:p ' + $wNs + '>' + $pPrHelv +
    (RunXml 'Note that the results for IA, ID, ME, and ND are based on too few samples to be significant. ' $true) +
    (RunXml 'E' $false) +
    (RunXml 'xcluding those states, there appears to be some difference in default rates by state that may be worth exploring.' $false) +
    (RunXml ' ' $true) +
    '</w:p>'
$p9.InsertXML($xml9)

# Paragraph 8: the inline OLE chart object - resize + refresh the
# (cosmetic) ObjectID in place, preserving the existing image/oleObject
# relationships untouched, and give the wrapping paragraph/run the same
# plain-Helvetica formatting as everywhere else.
$p8 = $d.Paragraphs.Item(8).Range
$xml8 = $p8.WordOpenXML
$xml8 = $xml8 -replace '(<w:body><w:p)\b[^>]*(>)<w:r>(<w:object)', ('$1$2' + $pPrHelv + '<w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/></w:rPr>$3')
$xml8 = $xml8.Replace('width:427pt;height:296pt', 'width:454.5pt;height:314.5pt')
$xml8 = $xml8.Replace('_1510307602', '_1510307834')
$p8.InsertXML($xml8)

# Paragraph 7: "Next, I analyzed ... "
$p7 = $d.Paragraphs.Item(7).Range
$xml7 = '<w:p ' + $wNs + '>' + $pPrHelv +
    (RunXml 'Next, I analyzed the percent of loans that are in default against those that are current on a state by state basis. The results are summarize below:' $false) +
    '</w:p>'
$p7.InsertXML($xml7)

# Paragraph 6: the FICO score paragraph
$p6 = $d.Paragraphs.Item(6).Range
$apos6 = [char]0x2019
$run6c = '. This suggests that there is no significant difference in default rates based on the FICO score of applicants who receive loans based on the Lending Club' + $apos6 + 's current selection criteria.'
$xml6 = '<w:p ' + $wNs + '>' + $pPrHelv +
    (RunXml 'The first, and somewhat obvious, dimension to consider was the average FICO score of the loan recipient. The average FICO score of recipients that are in default is 692 while the average score of those who are current is ' $true) +
    (RunXml '697' $false) +
    (RunXml $run6c $false) +
    '</w:p>'
$p6.InsertXML($xml6)

# Paragraph 5: empty spacer paragraph - remove entirely.
$d.Paragraphs.Item(5).Range.Delete()

# Paragraph 4: "I am exploring the Lending Club's loan dataset ... "
$p4 = $d.Paragraphs.Item(4).Range
$apos = [char]0x2019
$run4a = 'I am exploring the Lending Club' + $apos + 's loan dataset with the objective of looking for any trends in'
$xml4 = '<w:p ' + $wNs + '>' + $pPrHelv +
    (RunXml $run4a $false) +
    (RunXml ' loan' $true) +
    (RunXml ' defaults in order to create ' $true) +
    (RunXml 'a' $false) +
    (RunXml ' machine learning model to predict loan defaults. As a first step, I looked at some basic statistics and found several avenues to explore.' $true) +
    '</w:p>'
$p4.InsertXML($xml4)

# Paragraph 3: empty spacer paragraph - remove entirely.
$d.Paragraphs.Item(3).Range.Delete()

# Paragraph 2: "By now you've asked a bunch of questions ... " - remove
# entirely (its content is superseded by the new title paragraph below).
$d.Paragraphs.Item(2).Range.Delete()

# Paragraph 1: "Having made these plots ... " -> becomes the new title.
$p1 = $d.Paragraphs.Item(1).Range
$xml1 = '<w:p ' + $wNs + '>' + $pPrHelv +
    (RunXml 'The Lending Club: Initial Research' $false) +
    '</w:p>'
$p1.InsertXML($xml1)

Write-Output "edit applied"
